$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 31   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/25/2024  Through  3/31/2024"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "***.*"
$ws.Range("E15").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -31.25
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = 33.333333333333
$ws.Range("L16").Value = -12
$ws.Range("M16").Value = -34.328358208955
$ws.Range("N16").Value = -85.382059800664
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -62.5
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -52.777777777777
$ws.Range("I17").Value = 71
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = -26.80412371134
$ws.Range("L17").Value = -5.333333333333
$ws.Range("M17").Value = 36.538461538461
$ws.Range("N17").Value = -5.333333333333
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("D16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = 50
$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 23
$ws.Range("J18").Value = 27
$ws.Range("K18").Value = -14.814814814814
$ws.Range("L18").Value = -20.689655172413
$ws.Range("M18").Value = -65.671641791044
$ws.Range("N18").Value = -93.072289156626
$ws.Range("C19").Value = 3
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 35.294117647058
$ws.Range("I19").Value = 90
$ws.Range("J19").Value = 77
$ws.Range("K19").Value = 16.883116883116
$ws.Range("L19").Value = -11.764705882352
$ws.Range("M19").Value = 21.621621621621
$ws.Range("N19").Value = -29.6875
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 11.764705882352
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 67
$ws.Range("K20").Value = -19.402985074626
$ws.Range("L20").Value = -8.474576271186
$ws.Range("M20").Value = -11.475409836065
$ws.Range("N20").Value = -94.046306504961
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -43.478260869565
$ws.Range("F21").Value = 76
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = -19.148936170212
$ws.Range("I21").Value = 287
$ws.Range("J21").Value = 306
$ws.Range("K21").Value = -6.209150326797
$ws.Range("L21").Value = -9.748427672955
$ws.Range("M21").Value = -10.869565217391
$ws.Range("N21").Value = -83.628066172276
$ws.Range("C22").Value = 1
$ws.Range("C17").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 3
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -57.142857142857
$ws.Range("M22").Value = -50
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 47
$ws.Range("E24").Value = -57.446808510638
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 129
$ws.Range("H24").Value = -16.279069767441
$ws.Range("I24").Value = 345
$ws.Range("J24").Value = 381
$ws.Range("K24").Value = -9.448818897637
$ws.Range("L24").Value = 4.545454545454
$ws.Range("M24").Value = 75.126903553299
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 54.054054054054
$ws.Range("I25").Value = 169
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = 39.669421487603
$ws.Range("L25").Value = 32.03125
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 8.333333333333
$ws.Range("F26").Value = 53
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = 15.217391304347
$ws.Range("I26").Value = 146
$ws.Range("J26").Value = 131
$ws.Range("K26").Value = 11.450381679389
$ws.Range("L26").Value = 25.862068965517
$ws.Range("M26").Value = 2.097902097902
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = 0
$ws.Range("D28").Value = 3
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -83.333333333333
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = -43.75
$ws.Range("L28").Value = -30.76923076923
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E15").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("F31").Value = 1